$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '64.864.57'
$ws.Range('E2').Value = '  +0.21%  '

# Row 3
$ws.Range('D3').Value = '3.471.41'
$ws.Range('E3').Value = '  +0.47%  '

# Row 4
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.63'
$ws.Range('E5').Value = '  +0.09%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.94'

# Row 7
$ws.Range('E7').Value = '  +0.04%  '

# Row 8
$ws.Range('D8').Value = '3.470.82'
$ws.Range('E8').Value = '  +0.37%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.583'
$ws.Range('E9').Value = '  -6.52%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.23'
$ws.Range('E10').Value = '  +0.04%  '

# Row 11
$ws.Range('E11').Value = '  -2.22%  '

# Row 12
$ws.Range('E12').Value = '  -1.69%  '

# Row 13
$ws.Range('D13').Value = '4.067.24'
$ws.Range('E13').Value = '  +0.36%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.135'
$ws.Range('E14').Value = '  +0.12%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.66'
$ws.Range('E15').Value = '  -1.90%  '

# Row 16
$ws.Range('D16').Value = '64.928.97'
$ws.Range('E16').Value = '  +0.20%  '

# Row 17
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.471.82'
$ws.Range('E17').Value = '  -0.05%  '

# Row 18
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000169'
$ws.Range('E18').Value = '  -12.27%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.23'
$ws.Range('E19').Value = '  -3.37%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.86'
$ws.Range('E20').Value = '  -3.34%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '381.63'
$ws.Range('E21').Value = '  +0.32%  '

# Row 22
$ws.Range('E22').Value = '  -1.43%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.67'
$ws.Range('E23').Value = '  -0.01%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.997'
$ws.Range('E24').Value = '  -0.40%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.532'
$ws.Range('E25').Value = '  -3.55%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000121'
$ws.Range('E26').Value = '  +1.52%  '

# Row 27
$ws.Range('E27').Value = '  -0.41%  '

# Row 28
$ws.Range('E28').Value = '  +1.19%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.20%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.18'
$ws.Range('E30').Value = '  +1.07%  '

# Row 31
$ws.Range('E31').Value = '  -3.48%  '

# Row 32
$ws.Range('E32').Value = '  -1.09%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.31'
$ws.Range('E33').Value = '  -1.31%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.07'
$ws.Range('E34').Value = '  -1.37%  '

# Row 35
$ws.Range('E35').Value = '  -0.89%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '161.18'
$ws.Range('E36').Value = '  -0.02%  '

# Row 37
$ws.Range('E37').Value = '  -2.06%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0755'
$ws.Range('E38').Value = '  -2.66%  '

# Row 39
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '26.83'
$ws.Range('E39').Value = '  +0.88%  '

# Row 40
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '2.882.13'
$ws.Range('E40').Value = '  -2.40%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.815'
$ws.Range('E41').Value = '  +5.23%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.59'
$ws.Range('E42').Value = '  -0.19%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.53'
$ws.Range('E43').Value = '  -1.10%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.07'
$ws.Range('E44').Value = '  +0.79%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '25.98'
$ws.Range('E45').Value = '  +1.27%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0310'
$ws.Range('E46').Value = '  -2.37%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.47'
$ws.Range('E47').Value = '  +13.59%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '324.26'
$ws.Range('E48').Value = '  +5.10%  '

# Row 49
$ws.Range('E49').Value = '  -1.95%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.849'
$ws.Range('E50').Value = '  -2.19%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.48'
$ws.Range('E51').Value = '  -2.13%  '
